$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

$ws.Range("D7").Value = "Complete"
$ws.Range("F7").Value = "Trevor"

$ws.Range("B20").Value = 18
$ws.Range("C20").Value = "Recipe Unit Conversion"
$ws.Range("E20").Value = "Braden/Trevor"
$ws.Range("G20").Value = "Dependent"
$ws.Range("H20").Value = 6
$ws.Range("I20").Value = "Not Critical"

$ws.Range("H20").Select()
